# [Fonds de solidarite] Add 2020-08-12 data
# Update "nombre_aides" (column C) and "montant_total" (column D) figures
# for several region / legal-category rows to reflect the 2020-08-12 data
# refresh. Values are written as text (NumberFormat "@") so the exact
# original string representation (e.g. trailing ".00") is preserved,
# matching how the source data is stored in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue "C6" "543"
Set-TextValue "D6" "1576513.59"

Set-TextValue "C9" "166"
Set-TextValue "D9" "462230.44"

Set-TextValue "C10" "342"
Set-TextValue "D10" "1220255.08"

Set-TextValue "C11" "149"
Set-TextValue "D11" "496116.40"

Set-TextValue "C12" "6"
Set-TextValue "D12" "14450.00"

Set-TextValue "C14" "210"
Set-TextValue "D14" "554362.00"

Set-TextValue "C16" "472"
Set-TextValue "D16" "1696154.75"

Set-TextValue "C28" "258"
Set-TextValue "D28" "656542.64"

Set-TextValue "C30" "516"
Set-TextValue "D30" "2049800.70"

Set-TextValue "C32" "365"
Set-TextValue "D32" "1207198.17"

Set-TextValue "C33" "12"
Set-TextValue "D33" "38000.00"

Set-TextValue "C45" "341"
Set-TextValue "D45" "934867.74"

Set-TextValue "C47" "581"
Set-TextValue "D47" "2245831.99"

Set-TextValue "C48" "392"
Set-TextValue "D48" "1306867.16"

Set-TextValue "C91" "554"
Set-TextValue "D91" "1349722.15"

Set-TextValue "C93" "1061"
Set-TextValue "D93" "3537580.92"

Set-TextValue "C95" "966"
Set-TextValue "D95" "2898046.31"

Set-TextValue "C97" "46"
Set-TextValue "D97" "171183.00"
